$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32: G=5484
$ws.Range("H32").Value = 932.75
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 977
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 977
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -1629
# Row 98: G=36237
$ws.Range("H98").Value = 1763.762
$ws.Range("I98").Value = 1568.6
$ws.Range("K98").Value = 1568.6
$ws.Range("M98").Value = -70.59999999999991
# Row 106: G=19903
$ws.Range("H106").Value = 1974
$ws.Range("I106").Value = 1764.8334
$ws.Range("J106").Value = 3229
$ws.Range("K106").Value = 1764.8334
$ws.Range("L106").Value = 3229
$ws.Range("M106").Value = -1133.8334
$ws.Range("N106").Value = -4491
# Row 107: G=27766
$ws.Range("H107").Value = 671
$ws.Range("I107").Value = 671
$ws.Range("K107").Value = 671
$ws.Range("M107").Value = 1249
# Row 115: G=27957
$ws.Range("H115").Value = 710.5714
$ws.Range("I115").Value = 329
$ws.Range("K115").Value = 987
$ws.Range("M115").Value = 580
# Row 122: G=36237
$ws.Range("H122").Value = 1763.762
$ws.Range("I122").Value = 1568.6
$ws.Range("K122").Value = 4705.799999999999
$ws.Range("M122").Value = -2255.799999999999
# Row 132: G=44049
$ws.Range("H132").Value = 1358.8823
$ws.Range("I132").Value = 1033.6666
$ws.Range("K132").Value = 3100.9998
$ws.Range("M132").Value = -570.9998000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45: G=27714
$ws.Range("H45").Value = 1999.6666
$ws.Range("I45").Value = 1999.6666
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1999.6666
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1622.6666
$ws.Range("N45").ClearContents()
# Row 46: G=3498
$ws.Range("H46").Value = 3735
$ws.Range("I46").Value = 3735
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3735
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -3416
$ws.Range("N46").ClearContents()
# Row 61: G=43999
$ws.Range("H61").Value = 1933.3334
$ws.Range("I61").Value = 1933.3334
$ws.Range("K61").Value = 1933.3334
$ws.Range("M61").Value = -1721.3334
# Row 88: G=12530
$ws.Range("H88").Value = 2261.1
$ws.Range("I88").Value = 685.8
$ws.Range("J88").Value = 2786.2
$ws.Range("K88").Value = 685.8
$ws.Range("L88").Value = 2786.2
$ws.Range("M88").Value = -279.8
$ws.Range("N88").Value = -3598.2
# Row 91: G=12530
$ws.Range("H91").Value = 2261.1
$ws.Range("I91").Value = 685.8
$ws.Range("J91").Value = 2786.2
$ws.Range("K91").Value = 685.8
$ws.Range("L91").Value = 2786.2
$ws.Range("M91").Value = 718.2
$ws.Range("N91").Value = -5594.2
# Row 109: G=25646
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 110: G=27708
$ws.Range("H110").Value = 3748.2856
$ws.Range("I110").Value = 619.25
$ws.Range("K110").Value = 619.25
$ws.Range("M110").Value = 1425.75
# Row 132: G=43997
$ws.Range("H132").Value = 2598.1538
$ws.Range("I132").Value = 2598.1538
$ws.Range("K132").Value = 7794.4614
$ws.Range("M132").Value = -5264.4614
# Row 136: G=43999
$ws.Range("H136").Value = 1933.3334
$ws.Range("I136").Value = 1933.3334
$ws.Range("K136").Value = 5800.0002
$ws.Range("M136").Value = -3250.0002
# Row 139: G=42321
$ws.Range("H139").Value = 49999
$ws.Range("J139").Value = 49999
$ws.Range("L139").Value = 49999
$ws.Range("N139").Value = -60279

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94: G=19939
$ws.Range("H94").Value = 1166.4231
$ws.Range("I94").Value = 1240.3889
$ws.Range("K94").Value = 1240.3889
$ws.Range("M94").Value = -789.3888999999999
# Row 122: G=34096
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 134: G=43998
$ws.Range("H134").Value = 8814.177
$ws.Range("I134").Value = 9060.143
$ws.Range("K134").Value = 27180.429
$ws.Range("M134").Value = -24645.429

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58: G=44021
$ws.Range("H58").Value = 2052.6667
$ws.Range("I58").Value = 1979.1111
$ws.Range("K58").Value = 1979.1111
$ws.Range("M58").Value = -1776.1111
# Row 132: G=44019
$ws.Range("H132").Value = 1214.4615
$ws.Range("I132").Value = 1127.4
$ws.Range("J132").Value = 1504.6666
$ws.Range("K132").Value = 3382.2
$ws.Range("L132").Value = 4513.9998
$ws.Range("M132").Value = -852.2000000000003
$ws.Range("N132").Value = -9573.9998
# Row 136: G=44021
$ws.Range("H136").Value = 2052.6667
$ws.Range("I136").Value = 1979.1111
$ws.Range("K136").Value = 5937.3333
$ws.Range("M136").Value = -3387.3333

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 57: G=4655
$ws.Range("H57").Value = 5826.6665
$ws.Range("I57").Value = 4990
$ws.Range("J57").Value = 7500
$ws.Range("K57").Value = 14970
$ws.Range("L57").Value = 22500
$ws.Range("M57").Value = -14411
$ws.Range("N57").Value = -23618
# Row 62: G=12845
$ws.Range("H62").Value = 12500
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46372
# Row 65: G=12845
$ws.Range("H65").Value = 12500
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141864
# Row 98: G=19843
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()
# Row 109: G=27854
$ws.Range("H109").Value = 1742.6
$ws.Range("I109").Value = 1585.2
$ws.Range("K109").Value = 4755.6
$ws.Range("M109").Value = -3715.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 107: G=27802
$ws.Range("H107").Value = 383.83334
$ws.Range("I107").Value = 245.66667
$ws.Range("K107").Value = 245.66667
$ws.Range("M107").Value = 1674.33333
# Row 122: G=36182
$ws.Range("H122").Value = 6583295
$ws.Range("I122").Value = 7357388.5
$ws.Range("K122").Value = 22072165.5
$ws.Range("M122").Value = -22069715.5
# Row 126: G=36184
$ws.Range("H126").Value = 2999.75
$ws.Range("I126").Value = 2999.75
$ws.Range("K126").Value = 8999.25
$ws.Range("M126").Value = -6529.25
# Row 132: G=44008
$ws.Range("H132").Value = 1678.1111
$ws.Range("I132").Value = 1678.1111
$ws.Range("K132").Value = 5034.3333
$ws.Range("M132").Value = -2504.3333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40: G=36248
$ws.Range("H40").Value = 8335.5
$ws.Range("I40").Value = 8097.7144
$ws.Range("K40").Value = 8097.7144
$ws.Range("M40").Value = -7961.7144
# Row 125: G=34271
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
# Row 136: G=44060
$ws.Range("H136").Value = 3275.818
$ws.Range("I136").Value = 3275.818
$ws.Range("K136").Value = 9827.454000000002
$ws.Range("M136").Value = -7277.454000000002

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 5: G=3515
$ws.Range("H5").Value = 167266.67
$ws.Range("I5").Value = 167266.67
$ws.Range("K5").Value = 167266.67
$ws.Range("M5").Value = -167154.67
# Row 41: G=21725
$ws.Range("H41").Value = 15753
$ws.Range("I41").Value = 15918.2
$ws.Range("J41").Value = 15587.8
$ws.Range("K41").Value = 15918.2
$ws.Range("L41").Value = 15587.8
$ws.Range("M41").Value = -15528.2
$ws.Range("N41").Value = -16367.8
# Row 126: G=36210
$ws.Range("H126").Value = 2629.4119
$ws.Range("I126").Value = 2192.3076
$ws.Range("K126").Value = 6576.9228
$ws.Range("M126").Value = -4106.9228
